# Updates cryptos list price (column D) and 1h volume-change percentage
# (column E) for rows 2-51 of the active worksheet, matching the latest
# scrape from GitHub Actions.
#
# Column D prices are stored as plain text (not numbers) in the workbook,
# so each new value is written with a leading apostrophe (Excel's
# text-entry prefix) to force text interpretation, then the cell style is
# reset to "Normal" so no stray quote-prefix formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.181.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.61%  "
$ws.Range("D3").Value = "'2.384.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.00%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'550.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'135.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.12%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.26%  "
$ws.Range("D9").Value = "'2.382.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").Value = "'5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").Value = "'24.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").Value = "'2.811.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").Value = "'61.045.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.37%  "
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "'2.431.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.57%  "
$ws.Range("D19").Value = "'10.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'6.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.92%  "
$ws.Range("D22").Value = "'321.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'63.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "'0.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.55%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'8.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("E28").Value = "  +4.53%  "
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("D30").Value = "'0.0₃0757"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.97%  "
$ws.Range("D31").Value = "'171.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  +6.57%  "
$ws.Range("D33").Value = "'5.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("D34").Value = "'1.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.63%  "
$ws.Range("D35").Value = "'0.389"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.67%  "
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'4.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.62%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'330.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.80%  "
$ws.Range("E41").Value = "  +6.63%  "
$ws.Range("D42").Value = "'38.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "'146.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").Value = "'19.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.94%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'11.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("E51").Value = "  +5.26%  "
